# Apply EU-2024-develop branch edits to About and SoCiIEPTtB sheets
$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")
$soc = $wb.Worksheets.Item("SoCiIEPTtB")

# --- Rebuild the About sheet content ---
$about.Cells.Clear()

$about.Range("A1").Value = 'SoCiIEPTtB Share of Change in Industry Expenses Passed Through to Buyers'
$about.Range("A1").Font.Bold = $true
$about.Range("A3").Value = 'Source:'
$about.Range("A3").Font.Bold = $true
$about.Range("B3").Value = 'Anna Milanez, OECD'
$about.Range("B4").Value = 2017
$about.Range("B4").HorizontalAlignment = -4131
$about.Range("B5").Value = 'Legal tax liability, legal remittance responsibility and tax incidence: Three dimensions of business taxation'
$about.Range("B6").Value = 'https://www.oecd-ilibrary.org/deliver/e7ced3ea-en.pdf?itemId=%2Fcontent%2Fpaper%2Fe7ced3ea-en&mimeType=pdf'
$about.Range("B7").Value = 'Pages 40-41'
$about.Range("B9").Value = 'Cludius et al.'
$about.Range("B10").Value = 2020
$about.Range("B10").HorizontalAlignment = -4131
$about.Range("B11").Value = 'Ex-post investigation of cost pass-through in the EU ETS - an analysis for six industry sectors'
$about.Range("B12").Value = 'https://doi.org/10.1016/j.eneco.2020.104883'
$about.Range("A14").Value = 'Notes'
$about.Range("A14").Font.Bold = $true
$about.Range("A15").Value = 'This variable contains data on "tax incidence" - that is, who ultimately pays the burden of a tax, as opposed'
$about.Range("A16").Value = 'to who is legally responsible for physically remitting the tax.  In this case, it specifies the share'
$about.Range("A17").Value = 'of a tax that a business may pass on to consumers via increasing the prices of its products.'
$about.Range("A19").Value = 'Though the EPS applies this share to all policy-driven changes in business expenses, the data here come'
$about.Range("A20").Value = 'from sales taxes or VAT (value-added taxes), because this is what is available in the data, and money'
$about.Range("A21").Value = 'is fungible, so to the degree that businesses can pass on sales taxes, they should be able to pass on other'
$about.Range("A22").Value = 'policy-driven cost increases to the same degree.'
$about.Range("A24").Value = 'Most tax incidence data looks at corporate income taxes and how much of these taxes are passed on to'
$about.Range("A25").Value = 'a corporation''s workers (labor).  We use data on the share of taxes that are passed on to consumers,'
$about.Range("A26").Value = 'which is somewhat harder to find.'
$about.Range("A28").Value = 'The OECD paper cited above (by Anna Milanez) is the most comprehensive review I''ve seen, considering'
$about.Range("A29").Value = 'a number of countries and tax incidence types, including sales tax/VAT.  The key conclusions are presented'
$about.Range("A30").Value = 'as a written summary of other studies, rather than as a table of values.  Essentially, the key points are:'
$about.Range("B32").Value = 'Studies vary greatly in sales tax/VAT incidence estimates.  Several studies found that businesses can pass on'
$about.Range("B33").Value = 'more than 100% of tax increases due to market power.  One study found that in the short term,'
$about.Range("B34").Value = 'businesses could pass on about two thirds of the tax, but in the longer term, they were able to pass on'
$about.Range("B35").Value = 'all of the tax.'
$about.Range("B37").Value = 'The overall conclusion of a recent, multi-country study, "for changes in VAT standard rates, full'
$about.Range("B38").Value = 'pass-through (i.e., full incidence of indirect taxation on consumers) is a reasonable estimate"'
$about.Range("B39").Value = 'Benedek, de Mooij, Keen, & Wingender (2015)'
$about.Range("A41").Value = 'This is difficult to estimate, because it depends on minutiae of the tax design, and likely also the tax base'
$about.Range("A42").Value = '(e.g. whether or not there are substitute goods that escape the new tax), and other factors specific to'
$about.Range("A43").Value = 'geographies or industries.'
$about.Range("A45").Value = 'Accordingly, we adopt the recommendation of the Benedek, de Mooij, Keen, & Wingender (2015) study'
$about.Range("A46").Value = 'and assume full passthrough for many industries.'
$about.Range("A47").Value = 1
$about.Range("B47").Value = 'Share of cost increases passed through to consumers'
$about.Range("A49").Value = 'In the EU, however, we have found ETS to greatly impact prices, and therefore demand, imports, and exports. '
$about.Range("A50").Value = 'As a result, we use estimated values from Cludius et al. for some industries. '

# --- Rebuild the SoCiIEPTtB sheet content ---
$soc.Cells.Clear()

$soc.Range("A1").Value = 'Share of cost increases passed through to consumers'
$soc.Range("A1").Font.Italic = $true
$soc.Range("B1").Value = 'Unit: dimensionless (% passthrough)'

$soc.Range("A2").Value = 'agriculture and forestry 01T03'
$soc.Range("B2").Formula = "=About!A`$47"
$soc.Range("A3").Value = 'coal mining 05'
$soc.Range("B3").Formula = "=About!A`$47"
$soc.Range("A4").Value = 'oil and gas extraction 06'
$soc.Range("B4").Formula = "=About!A`$47"
$soc.Range("A5").Value = 'other mining and quarrying 07T08'
$soc.Range("B5").Formula = "=About!A`$47"
$soc.Range("A6").Value = 'food beverage and tobacco 10T12'
$soc.Range("B6").Formula = "=About!A`$47"
$soc.Range("A7").Value = 'textiles apparel and leather 13T15'
$soc.Range("B7").Formula = "=About!A`$47"
$soc.Range("A8").Value = 'wood products 16'
$soc.Range("B8").Formula = "=About!A`$47"
$soc.Range("A9").Value = 'pulp paper and printing 17T18'
$soc.Range("B9").Formula = "=About!A`$47"
$soc.Range("A10").Value = 'refined petroleum and coke 19'
$soc.Range("B10").Formula = "=About!A`$47"
$soc.Range("A11").Value = 'chemicals 20'
$soc.Range("B11").Formula = "=About!A`$47"
$soc.Range("A12").Value = 'rubber and plastic products 22'
$soc.Range("B12").Formula = "=About!A`$47"
$soc.Range("A13").Value = 'glass and glass products 231'
$soc.Range("B13").Value = 0.5
$soc.Range("A14").Value = 'cement and other nonmetallic minerals 239'
$soc.Range("B14").Value = 0.3
$soc.Range("A15").Value = 'iron and steel 241'
$soc.Range("B15").Value = 0.7
$soc.Range("A16").Value = 'other metals 242'
$soc.Range("B16").Formula = "=About!A`$47"
$soc.Range("A17").Value = 'metal products except machinery and vehicles 25'
$soc.Range("B17").Formula = "=About!A`$47"
$soc.Range("A18").Value = 'computers and electronics 26'
$soc.Range("B18").Formula = "=About!A`$47"
$soc.Range("A19").Value = 'appliances and electrical equipment 27'
$soc.Range("B19").Formula = "=About!A`$47"
$soc.Range("A20").Value = 'other machinery 28'
$soc.Range("B20").Formula = "=About!A`$47"
$soc.Range("A21").Value = 'road vehicles 29'
$soc.Range("B21").Formula = "=About!A`$47"
$soc.Range("A22").Value = 'nonroad vehicles 30'
$soc.Range("B22").Formula = "=About!A`$47"
$soc.Range("A23").Value = 'other manufacturing 31T33'
$soc.Range("B23").Formula = "=About!A`$47"
$soc.Range("A24").Value = 'energy pipelines and gas processing 352T353'
$soc.Range("B24").Formula = "=About!A`$47"
$soc.Range("A25").Value = 'water and waste 36T39'
$soc.Range("B25").Formula = "=About!A`$47"
$soc.Range("A26").Value = 'construction 41T43'
$soc.Range("B26").Formula = "=About!A`$47"

# --- Selection / active sheet state ---
$about.Range("A51").Select()
$soc.Range("B14").Select()
$soc.Activate()

$wb.Application.Calculate()
